$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date stamp for every
# data row (rows 2-306). The whole column is bumped from 2023-09-19
# (serial 45188) to 2023-09-20 (serial 45189).
$lastRow = 306

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45189
}
